$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new progress entry for Tejomay Padole (row 7) under the "21-Feb" column (D)
$ws.Range("D7").Value = "Studied firebase and how to use it with flutter"

# Reflect the new active cell selection on Sheet1
$ws.Range("D7").Select()
